# Sync attendance_reports: normalize the "Recorded By" (column G) author
# ordering for specific known value combinations.
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"             -> "admin@admin.com, System"
#   "backup@backdoor.com, System, system" -> "system, backup@backdoor.com, System"
# All other "Recorded By" values are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$map = @{
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "backup@backdoor.com, System, system" = "system, backup@backdoor.com, System";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }

    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
